$d = $word.ActiveDocument

# Locate the "Unintentional Drowning Data Brief" hyperlink inside the
# Projects/Experience bullet: "... data products (Unintentional Drowning
# Data Brief Link) for dissemination to key stakeholders."
$h = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.TextToDisplay -eq "Unintentional Drowning Data Brief") {
        $h = $candidate
        break
    }
}
if ($h -eq $null) {
    $h = $d.Hyperlinks.Item(1)
}

$hs = $h.Range.Start
$he = $h.Range.End

# Find where the trailing " Link) " text (right after the hyperlink) ends,
# i.e. where "for dissemination" starts, so we can remove exactly that span.
$searchRng = $d.Range($he, $d.Content.End)
$searchRng.Find.Execute("for dissemination", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterEnd = $searchRng.Start

# 1) Remove the " Link) " text trailing the hyperlink.
$rngAfter = $d.Range($he, $afterEnd)
$rngAfter.Text = ""

# 2) Remove the hyperlinked run "Unintentional Drowning Data Brief" itself
#    (this also drops the now-empty w:hyperlink wrapper / its relationship
#    reference from the run content).
$rngLink = $d.Range($hs, $he)
$rngLink.Text = ""

# 3) The run right before the (now removed) hyperlink held " (" -- collapse
#    it down to a single space so the sentence reads "...data products for
#    dissemination...".
$rngBefore = $d.Range($hs - 2, $hs)
$rngBefore.Text = " "
